$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.336.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.677.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.06%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5103"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06319"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07357"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.684.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.532"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5764"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.908.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.88%  "

$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -13.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.375.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.72%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("E21").Value = "  -4.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.228"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.480"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.57%  "

$ws.Range("E27").Value = "  -6.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.336"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05842"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.331"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.87%  "

$ws.Range("E32").Value = "  -6.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.500"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.653"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5943"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.357"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.82%  "

$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.094.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.892"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8590"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.832.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.54%  "

$ws.Range("E46").Value = "  +4.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.019"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4311"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.98%  "
